$wb = $excel.ActiveWorkbook

# The "natural gas nonpeaker" plant type is being split into two separate
# plant types: "natural gas steam turbine" and "natural gas combined cycle".
# This affects both the BPaFF-BITPTaP and BPaFF-BDTPTPF sheets identically.

foreach ($sheetName in @("BPaFF-BITPTaP", "BPaFF-BDTPTPF")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3 currently holds "natural gas nonpeaker" (value 0).
    # Insert a new row above it so we end up with two rows:
    #   row 3: natural gas steam turbine = 0
    #   row 4: natural gas combined cycle = 0 (was "natural gas nonpeaker")
    $ws.Rows.Item(3).Insert()
    $ws.Range("A3:B3").ClearFormats()

    $ws.Range("A3").Value = "natural gas steam turbine"
    $ws.Range("B3").Value = 0

    $ws.Range("A4").Value = "natural gas combined cycle"
    $ws.Range("B4").Value = 0
}
